$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.028, -0.07099999999999999, 0.844, 85, -6, 122),
    @(-0.046, -0.363, -0.414, -106, -14, 119),
    @(0.202, 0.11, -0.159, -95, -68, -139),
    @(0.134, 0.126, 0.377, 22, 17, 6),
    @(0.673, -0.428, 0.285, -162, -51, -126),
    @(0.171, -0.022, -0.66, 94, -47, 97),
    @(0.008, 0.478, 0.544, -9, -22, -95),
    @(0.09, 0.128, 0.369, 18, -38, -171),
    @(0.304, -0.309, -0.37, 140, -27, 56),
    @(0.344, 0.093, 0.6860000000000001, 77, 18, 139)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

$wb.Save()
